# Applies crypto price/volume updates per commit "Updated cryptos list on Fri Sep  6 23:49:19 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must remain plain TEXT even when it looks numeric
# (Excel normally auto-converts numeric-looking strings into numbers; forcing the
# cell to Text format before the assignment keeps it literal, then the style is
# reset back to Normal so no stray formatting is left behind).
function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '53.614.01'
$ws.Range("E2").Value = '  -4.43%  '
$ws.Range("D3").Value = '2.198.23'
$ws.Range("E3").Value = '  -6.99%  '
Set-TextValue $ws.Range("D4") '1.00'
$ws.Range("E4").Value = '  +0.47%  '
Set-TextValue $ws.Range("D5") '486.43'
$ws.Range("E5").Value = '  -3.24%  '
Set-TextValue $ws.Range("D6") '124.85'
$ws.Range("E6").Value = '  -3.63%  '
Set-TextValue $ws.Range("D7") '1.00'
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  -4.27%  '
$ws.Range("D9").Value = '2.219.62'
$ws.Range("E9").Value = '  -6.18%  '
Set-TextValue $ws.Range("D10") '0.0924'
$ws.Range("E10").Value = '  -6.24%  '
$ws.Range("E12").Value = '  -3.31%  '
$ws.Range("E13").Value = '  -3.30%  '
$ws.Range("D14").Value = '2.593.58'
$ws.Range("E14").Value = '  -6.81%  '
Set-TextValue $ws.Range("D15") '21.10'
$ws.Range("E15").Value = '  -1.36%  '
$ws.Range("D16").Value = '53.603.82'
$ws.Range("E16").Value = '  -4.35%  '
$ws.Range("E17").Value = '  -3.19%  '
$ws.Range("D18").Value = '2.227.40'
$ws.Range("E18").Value = '  -6.54%  '
$ws.Range("E19").Value = '  -1.56%  '
$ws.Range("E20").Value = '  -4.51%  '
Set-TextValue $ws.Range("D21") '294.91'
$ws.Range("E21").Value = '  -3.99%  '
Set-TextValue $ws.Range("D22") '6.16'
$ws.Range("E22").Value = '  -2.36%  '
$ws.Range("E23").Value = '  +0.08%  '
Set-TextValue $ws.Range("D24") '62.92'
$ws.Range("E24").Value = '  -4.66%  '
Set-TextValue $ws.Range("D25") '0.997'
$ws.Range("E25").Value = '  -0.41%  '
Set-TextValue $ws.Range("D26") '0.366'
$ws.Range("E26").Value = '  -0.85%  '
$ws.Range("D27").Value = '2.299.98'
$ws.Range("E27").Value = '  -6.98%  '
$ws.Range("E28").Value = '  -1.39%  '
Set-TextValue $ws.Range("D29") '7.02'
$ws.Range("E29").Value = '  -3.03%  '
Set-TextValue $ws.Range("D30") '165.12'
$ws.Range("E30").Value = '  -4.35%  '
$ws.Range("E31").Value = '  -3.96%  '
Set-TextValue $ws.Range("D32") '0.999'
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("B33").Value = 'PEPE'
$ws.Range("C33").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D33").Value = '0.0₃0663'
$ws.Range("E33").Value = '  -6.77%  '
$ws.Range("B34").Value = 'FirstDigitalUSD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range("D34") '0.994'
$ws.Range("E34").Value = '  -0.35%  '
Set-TextValue $ws.Range("D35") '5.70'
$ws.Range("E35").Value = '  -1.49%  '
$ws.Range("E36").Value = '  -1.26%  '
Set-TextValue $ws.Range("D37") '17.30'
$ws.Range("E37").Value = '  -1.90%  '
$ws.Range("E38").Value = '  -1.29%  '
Set-TextValue $ws.Range("D39") '0.832'
$ws.Range("E39").Value = '  +4.15%  '
$ws.Range("E40").Value = '  -4.68%  '
Set-TextValue $ws.Range("D41") '35.82'
$ws.Range("E41").Value = '  -1.25%  '
Set-TextValue $ws.Range("D42") '0.367'
$ws.Range("E42").Value = '  -0.63%  '
$ws.Range("E43").Value = '  -1.36%  '
Set-TextValue $ws.Range("D44") '126.47'
$ws.Range("E44").Value = '  -2.12%  '
$ws.Range("E45").Value = '  -2.63%  '
Set-TextValue $ws.Range("D46") '4.71'
$ws.Range("E46").Value = '  +0.53%  '
Set-TextValue $ws.Range("D47") '0.0880'
$ws.Range("E47").Value = '  -2.58%  '
Set-TextValue $ws.Range("D48") '0.535'
$ws.Range("E48").Value = '  -4.79%  '
Set-TextValue $ws.Range("D49") '231.60'
$ws.Range("E49").Value = '  -2.74%  '
$ws.Range("E50").Value = '  -2.25%  '
$ws.Range("E51").Value = '  -3.46%  '
